$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert six new rows before the old row 38 ("Bond" section), shifting
#    everything below down by six rows (38 -> 44, etc.).
# ---------------------------------------------------------------------------
$ws.Rows("38:43").Insert()

# ---------------------------------------------------------------------------
# 2. New section header row 38: "Alias"
#    A38 reuses the "Neutral" highlighted style that A1/B1 used to have;
#    B38:E38 reuse the plain fill-highlighted continuation style already
#    used elsewhere in the sheet (e.g. row 68, col B "DOB").
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("A38").Value = "Alias"

$ws.Range("B68").Copy()
$ws.Range("B38:E38").PasteSpecial(-4122)

$ws.Range("A38:E38").Font.Bold = $true
$ws.Range("A38:E38").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Data rows 39-43. Columns A-D reuse the common plain style (e.g. B10).
#    Column E reuses that same plain style first, then gets bolded with an
#    explicit font/fill so Excel mints the new "Alias answer" font+fill
#    combination actually used by the workbook.
# ---------------------------------------------------------------------------
$ws.Range("B10").Copy()
$ws.Range("A39:D43").PasteSpecial(-4122)

$ws.Range("B68").Copy()
$ws.Range("E39:E43").PasteSpecial(-4122)
$ws.Range("E39:E43").Font.Bold = $true
$ws.Range("E39:E43").Font.Name = "Calibri"
$ws.Range("E39:E43").Font.Size = 12

$ws.Range("A39:E43").RowHeight = 63

# ---------------------------------------------------------------------------
# 4. Cell values.
# ---------------------------------------------------------------------------
$ws.Range("B39").Value = "DOB"
$ws.Range("C39").Value = "Person DOB"
$ws.Range("E39").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Identity[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:PersonAliasIdentityAssociation/nc:Identity/@structures:ref]/nc:IdentityPersonRepresentation/nc:PersonBirthDate/nc:Date"

$ws.Range("B40").Value = "First name"
$ws.Range("C40").Value = "Person given name"
$ws.Range("E40").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Identity[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:PersonAliasIdentityAssociation/nc:Identity/@structures:ref]/nc:IdentityPersonRepresentation/nc:PersonName/nc:PersonGivenName"

$ws.Range("B41").Value = "Middle name"
$ws.Range("C41").Value = "Person middle name"
$ws.Range("E41").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Identity[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:PersonAliasIdentityAssociation/nc:Identity/@structures:ref]/nc:IdentityPersonRepresentation/nc:PersonName/nc:PersonMiddleName"

$ws.Range("B42").Value = "Last name"
$ws.Range("C42").Value = "Person last name"
$ws.Range("E42").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Identity[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:PersonAliasIdentityAssociation/nc:Identity/@structures:ref]/nc:IdentityPersonRepresentation/nc:PersonName/nc:PersonSurName"

$ws.Range("B43").Value = "Sex"
$ws.Range("C43").Value = "Person Sex"
$ws.Range("E43").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Identity[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:PersonAliasIdentityAssociation/nc:Identity/@structures:ref]/nc:IdentityPersonRepresentation/j:PersonSexCode"

# ---------------------------------------------------------------------------
# 5. The title row (A1:B1) loses its bold "Neutral" highlight (that look now
#    belongs to the new "Alias" section header) and becomes plain wrapped
#    text.
# ---------------------------------------------------------------------------
$ws.Range("A1:B1").Font.Bold = $false
$ws.Range("A1:B1").WrapText = $true

# ---------------------------------------------------------------------------
# 6. Update the frozen-pane view to where the user ended up working.
# ---------------------------------------------------------------------------
$ws.Range("A38").Select()
$excel.ActiveWindow.ScrollRow = 38
$ws.Range("A42").Select()

Write-Host "Alias section inserted"
